# Apply weekly update: insert a new weekly record row for
# "Macroferia Regional de Talca - Acelga" as row 132, shifting the
# existing rows 132:175 down to 133:176.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 132 (pushes old row 132.. down by one)
$ws.Rows.Item(132).Insert()

# Fill in the new row 132 with this week's data
$ws.Cells.Item(132, 1).Value  = 5
$ws.Cells.Item(132, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(132, 3).Value  = "Maule"
$ws.Cells.Item(132, 4).Value  = 44468
$ws.Cells.Item(132, 5).Value  = 7
$ws.Cells.Item(132, 6).Value  = 100112009
$ws.Cells.Item(132, 7).Value  = "Acelga"
$ws.Cells.Item(132, 8).Value  = "Sin especificar"
$ws.Cells.Item(132, 9).Value  = "Primera"
$ws.Cells.Item(132, 10).Value = 500
$ws.Cells.Item(132, 11).Value = 2000
$ws.Cells.Item(132, 12).Value = 2000
$ws.Cells.Item(132, 13).Value = 2000
$ws.Cells.Item(132, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(132, 15).Value = "Región del Maule"
$ws.Cells.Item(132, 16).Value = 500
$ws.Cells.Item(132, 17).Value = 4
$ws.Cells.Item(132, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date/time number format used
# throughout column D
$ws.Cells.Item(132, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
